$d = $word.ActiveDocument

# 1. Insert a new bold paragraph "Play Action Dragons for Free - Review & Features"
#    right before the final ("Prompt: ...") paragraph. First create a genuinely
#    new (empty) paragraph immediately before the last paragraph - this keeps all
#    the other paragraphs completely untouched - and then stamp that brand new,
#    isolated, empty paragraph with the exact OOXML run structure we need (a
#    leading empty run followed by the bold text run).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$null = $lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item($count)
$newParaRange = $newPara.Range.Duplicate
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Action Dragons for Free - Review &amp; Features</w:t></w:r></w:p>'
$null = $newParaRange.InsertXML($newParaXml)

# 2. Replace the text of the final ("Prompt: ...") paragraph's run with the
#    meta-description text, keeping its existing italic formatting intact.
$oldPromptText = 'Prompt: Create a feature image for "Action Dragons" that fits with the game''s theme and visually represents the unique features and symbols of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses. The image should include the dragon wild symbol and the green door symbol that triggers the Lucky Multiplier Bonus Spins. The background should be inspired by Chinese culture, with elements of jade and traditional Chinese music. The image should also feature symbols from the game, such as gold coins, vases, and perhaps the yin and yang symbol. The Maya warrior should be energetic and happy, holding a golden coin with the dragon wild symbol on it. The warrior should also be wearing glasses to represent a modern technology twist on the ancient theme. The image should be exciting and dynamic, capturing the game''s potential for high payouts and unique bonuses.'
$newPromptText = 'Discover Action Dragons slot game and its unique bonuses. Play free now and increase rewards by up to 5x.'
$null = $d.Content.Find.Execute($oldPromptText, $true, $false, $false, $false, $false, $true, 1, $false, $newPromptText, 2)

# 3. Remove the original "Meta description" paragraph (the second paragraph in
#    the document, right after the title heading).
$metaPara = $d.Paragraphs.Item(2)
$null = $metaPara.Range.Delete()
